$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112129065
$ws.Range("B2").Value = 96720
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "12"
$ws.Range("Q2").Value = 655221
$ws.Range("R2").Value = 6675131
$ws.Range("Z2").Value = "10:12"
$ws.Range("AB2").Value = "10:13"
$ws.Range("AC2").Value = "Djupt nere bland ris och mossa."

# Row 3
$ws.Range("A3").Value = 112129079
$ws.Range("B3").Value = 98961
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = "Blåsippa"
$ws.Range("G3").Value = "Hepatica nobilis"
$ws.Range("H3").Value = "Schreb."
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("Q3").Value = 655188
$ws.Range("R3").Value = 6675131
$ws.Range("Z3").Value = "10:28"
$ws.Range("AB3").Value = "10:28"

# Row 5
$ws.Range("B5").Value = 96720

# Row 6
$ws.Range("A6").Value = 112129067
$ws.Range("B6").Value = 96720
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "3"
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("Q6").Value = 655241
$ws.Range("R6").Value = 6675125
$ws.Range("Z6").Value = "10:16"
$ws.Range("AB6").Value = "10:17"
$ws.Range("AC6").Value = "Tuff tillvaro nära hyggeskanten."

# Row 7
$ws.Range("A7").Value = 112129069
$ws.Range("B7").Value = 96720
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "2"
$ws.Range("Q7").Value = 655168
$ws.Range("R7").Value = 6675142
$ws.Range("Z7").Value = "10:37"
$ws.Range("AB7").Value = "10:38"
$ws.Range("AC7").ClearContents()

# Row 8
$ws.Range("A8").Value = 112129072
$ws.Range("B8").Value = 96720
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "8"
$ws.Range("Q8").Value = 655162
$ws.Range("R8").Value = 6675144
$ws.Range("Z8").Value = "10:39"
$ws.Range("AB8").Value = "10:41"
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 112129063
$ws.Range("B9").Value = 96720
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "38"
$ws.Range("J9").Value = "plantor/tuvor"
$ws.Range("M9").ClearContents()
$ws.Range("Q9").Value = 655225
$ws.Range("R9").Value = 6675110
$ws.Range("Z9").Value = "10:06"
$ws.Range("AB9").Value = "10:07"
$ws.Range("AC9").Value = "Tätt med småplantor."

# Row 10
$ws.Range("A10").Value = 112128602
$ws.Range("B10").Value = 56575
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 103021
$ws.Range("F10").Value = "Talltita"
$ws.Range("G10").Value = "Poecile montanus"
$ws.Range("H10").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "1"
$ws.Range("J10").ClearContents()
$ws.Range("M10").Value = "permanent revir"
$ws.Range("Q10").Value = 655214
$ws.Range("R10").Value = 6675119
$ws.Range("Z10").Value = "10:09"
$ws.Range("AB10").Value = "10:09"
